$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1377.4814
$ws.Range("I19").Value = 710.0625
$ws.Range("J19").Value = 2348.2727
$ws.Range("K19").Value = 710.0625
$ws.Range("L19").Value = 2348.2727
$ws.Range("M19").Value = -535.0625
$ws.Range("N19").Value = -2698.2727
$ws.Range("H70").Value = 1400.2142
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2730
$ws.Range("H73").Value = 1400.2142
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -2064
$ws.Range("H80").Value = 2141
$ws.Range("I80").Value = 1794.7142
$ws.Range("J80").Value = 2343
$ws.Range("K80").Value = 5384.142599999999
$ws.Range("L80").Value = 7029
$ws.Range("M80").Value = -4386.142599999999
$ws.Range("N80").Value = -9025
$ws.Range("H83").Value = 2141
$ws.Range("I83").Value = 1794.7142
$ws.Range("J83").Value = 2343
$ws.Range("K83").Value = 16152.4278
$ws.Range("L83").Value = 21087
$ws.Range("M83").Value = -11160.4278
$ws.Range("N83").Value = -31071
$ws.Range("H88").Value = 6799.467
$ws.Range("J88").Value = 8430.888999999999
$ws.Range("L88").Value = 8430.888999999999
$ws.Range("N88").Value = -9242.888999999999
$ws.Range("H91").Value = 6799.467
$ws.Range("J91").Value = 8430.888999999999
$ws.Range("L91").Value = 8430.888999999999
$ws.Range("N91").Value = -11238.889
$ws.Range("H131").Value = 333335200
$ws.Range("I131").Value = 333335200
$ws.Range("K131").Value = 1000005600
$ws.Range("M131").Value = -1000000560
$ws.Range("H138").Value = 2152812.2
$ws.Range("I138").Value = 997.8570999999999
$ws.Range("J138").Value = 3451320.8
$ws.Range("K138").Value = 2993.5713
$ws.Range("L138").Value = 10353962.4
$ws.Range("M138").Value = 2146.4287
$ws.Range("N138").Value = -10364242.4
$ws.Range("H139").Value = 37779.668
$ws.Range("J139").Value = 37779.668
$ws.Range("L139").Value = 37779.668
$ws.Range("N139").Value = -48059.668

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1347.4681
$ws.Range("I86").Value = 1214.4
$ws.Range("K86").Value = 1214.4
$ws.Range("M86").Value = -91.40000000000009
$ws.Range("H89").Value = 1347.4681
$ws.Range("I89").Value = 1214.4
$ws.Range("K89").Value = 6072
$ws.Range("M89").Value = -456
$ws.Range("H134").Value = 432416.22
$ws.Range("I134").Value = 608125.6
$ws.Range("J134").Value = 2904.2593
$ws.Range("K134").Value = 1824376.8
$ws.Range("L134").Value = 8712.777900000001
$ws.Range("M134").Value = -1821841.8
$ws.Range("N134").Value = -13782.7779

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1464.6227
$ws.Range("I31").Value = 848.8484999999999
$ws.Range("K31").Value = 848.8484999999999
$ws.Range("M31").Value = -553.8484999999999
$ws.Range("H34").Value = 1464.6227
$ws.Range("I34").Value = 848.8484999999999
$ws.Range("K34").Value = 848.8484999999999
$ws.Range("M34").Value = -646.8484999999999
$ws.Range("H86").Value = 7593.9487
$ws.Range("I86").Value = 6437.1924
$ws.Range("J86").Value = 9907.462
$ws.Range("K86").Value = 6437.1924
$ws.Range("L86").Value = 9907.462
$ws.Range("M86").Value = -5314.1924
$ws.Range("N86").Value = -12153.462
$ws.Range("H89").Value = 7593.9487
$ws.Range("I89").Value = 6437.1924
$ws.Range("J89").Value = 9907.462
$ws.Range("K89").Value = 32185.962
$ws.Range("L89").Value = 49537.31
$ws.Range("M89").Value = -26569.962
$ws.Range("N89").Value = -60769.31
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H107").Value = 2185.8235
$ws.Range("I107").Value = 610
$ws.Range("J107").Value = 2842.4167
$ws.Range("K107").Value = 610
$ws.Range("L107").Value = 2842.4167
$ws.Range("M107").Value = 1310
$ws.Range("N107").Value = -6682.4167
$ws.Range("H132").Value = 447354.4
$ws.Range("I132").Value = 1173.7164
$ws.Range("J132").Value = 2315736
$ws.Range("K132").Value = 3521.1492
$ws.Range("L132").Value = 6947208
$ws.Range("M132").Value = -991.1491999999998
$ws.Range("N132").Value = -6952268

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1212613
$ws.Range("I113").Value = 2020700.4
$ws.Range("J113").Value = 482.1
$ws.Range("K113").Value = 6062101.199999999
$ws.Range("L113").Value = 1446.3
$ws.Range("M113").Value = -6059931.199999999
$ws.Range("N113").Value = -5786.3

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2551.3333
$ws.Range("I80").Value = 2341
$ws.Range("J80").Value = 3603
$ws.Range("K80").Value = 2341
$ws.Range("L80").Value = 3603
$ws.Range("M80").Value = -1343
$ws.Range("N80").Value = -5599
$ws.Range("H83").Value = 2551.3333
$ws.Range("I83").Value = 2341
$ws.Range("J83").Value = 3603
$ws.Range("K83").Value = 11705
$ws.Range("L83").Value = 18015
$ws.Range("M83").Value = -6713
$ws.Range("N83").Value = -27999
$ws.Range("H102").Value = 1612.4166
$ws.Range("I102").Value = 1290.5
$ws.Range("J102").Value = 2256.25
$ws.Range("K102").Value = 1290.5
$ws.Range("L102").Value = 2256.25
$ws.Range("M102").Value = 331.5
$ws.Range("N102").Value = -5500.25
$ws.Range("H132").Value = 1889035.5
$ws.Range("I132").Value = 2176.475
$ws.Range("J132").Value = 7694755.5
$ws.Range("K132").Value = 6529.424999999999
$ws.Range("L132").Value = 23084266.5
$ws.Range("M132").Value = -3999.424999999999
$ws.Range("N132").Value = -23089326.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 875.8946999999999
$ws.Range("I22").Value = 940
$ws.Range("K22").Value = 940
$ws.Range("M22").Value = -645
$ws.Range("H27").Value = 875.8946999999999
$ws.Range("I27").Value = 940
$ws.Range("K27").Value = 940
$ws.Range("M27").Value = -833
$ws.Range("H55").Value = 509.70834
$ws.Range("I55").Value = 208.22223
$ws.Range("J55").Value = 690.6
$ws.Range("K55").Value = 208.22223
$ws.Range("L55").Value = 690.6
$ws.Range("M55").Value = -35.22223
$ws.Range("N55").Value = -1036.6
$ws.Range("H101").Value = 16685.25
$ws.Range("J101").Value = 16685.25
$ws.Range("L101").Value = 16685.25
$ws.Range("N101").Value = -23175.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 28457.572
$ws.Range("J103").Value = 28457.572
$ws.Range("L103").Value = 28457.572
$ws.Range("N103").Value = -30801.572
